# Rename "CIN_waffle" -> "CIN_waffle_sv", and insert a new sheet
# "CIN_waffle_any" right after it (before "CSC_timing"), populated with
# the "any offence" waffle-chart source data.

$wb = $excel.ActiveWorkbook

# 1) Rename the existing waffle sheet to the "sv" (serious violence) variant.
$wsSv = $wb.Worksheets.Item("CIN_waffle")
$wsSv.Name = "CIN_waffle_sv"

# 2) Insert the new "any offence" waffle sheet directly after it.
$wsAny = $wb.Worksheets.Add($null, $wsSv)
$wsAny.Name = "CIN_waffle_any"

# Reuse the header row's formatting (bold + centered) from the sibling sheet
# so the new header row matches the existing style exactly.
$wsSv.Range("A1:F1").Copy()
$wsAny.Range("A1:F1").PasteSpecial(-4122)

# Header labels
$wsAny.Range("A1").Value = "LA"
$wsAny.Range("B1").Value = "indicator"
$wsAny.Range("C1").Value = "propany_count_CIN"
$wsAny.Range("D1").Value = "also_propany_count_CIN"
$wsAny.Range("E1").Value = "propany_count_not_CIN"
$wsAny.Range("F1").Value = "not_also_propany_count_CIN"

# Data rows: LA, indicator, propany_count_CIN, also_propany_count_CIN,
#            propany_count_not_CIN, not_also_propany_count_CIN
$rows = @(
  @("Birmingham","School",19,14,81,86),
  @("Bradford","School",20,17,80,83),
  @("Haringey","School",19,12,81,88),
  @("Lambeth","School",37,22,63,78),
  @("Leeds","School",35,20,65,80),
  @("Liverpool","School",34,12,66,88),
  @("Manchester","School",43,26,57,74),
  @("Newham","School",33,18,67,82),
  @("Sheffield","School",35,24,65,76),
  @("Southwark","School",16,13,84,87),
  @("Birmingham","Home",19,15,81,85),
  @("Bradford","Home",21,18,79,82),
  @("Haringey","Home",28,20,72,80),
  @("Lambeth","Home",32,21,68,79),
  @("Leeds","Home",36,21,64,79),
  @("Liverpool","Home",39,16,61,84),
  @("Manchester","Home",38,23,62,77),
  @("Newham","Home",33,18,67,82),
  @("Sheffield","Home",34,23,66,77),
  @("Southwark","Home",17,14,83,86)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $wsAny.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

Write-Host "CIN_waffle renamed to CIN_waffle_sv; CIN_waffle_any inserted with data."
